$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 12378
$ws.Range("C3:C4").Value = 11236
$ws.Range("C5:C17").Value = 10447
$ws.Range("C18:C20").Value = 10273
$ws.Range("C21").Value = 9857
$ws.Range("C22:C23").Value = 9575
$ws.Range("C24:C28").Value = 9544
$ws.Range("C29:C35").Value = 8820
$ws.Range("C36:C56").Value = 8336
$ws.Range("C57").Value = 8197
$ws.Range("C58:C83").Value = 7817
$ws.Range("C84:C115").Value = 7723
$ws.Range("C186:C208").Value = 7534
$ws.Range("C209:C252").Value = 7345
